$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 305, shifting existing row 305 (and below) down to 306.
$ws.Rows.Item(305).Insert()

# Populate the newly inserted row 305 with the new record's data.
$ws.Range("A305").Value = 4
$ws.Range("B305").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C305").Value = "Los Lagos"
$ws.Range("D305").Value = 44795
$ws.Range("E305").Value = 10
$ws.Range("F305").Value = 100114013
$ws.Range("G305").Value = "Zanahoria"
$ws.Range("H305").Value = "Sin especificar"
$ws.Range("I305").Value = "Primera"
$ws.Range("J305").Value = 300
$ws.Range("K305").Value = 9500
$ws.Range("L305").Value = 10000
$ws.Range("M305").Value = 9750
$ws.Range("N305").Value = "`$/saco 20 kilos"
$ws.Range("O305").Value = "Provincia de Llanquihue"
$ws.Range("P305").Value = 488
$ws.Range("Q305").Value = 20
$ws.Range("R305").Value = "Hortaliza"
